# Generate Report for Handoff
# Updates status from "In Translation" -> "Ready for handoff" and refreshes
# the "Latest Handoff"/"Latest HO Xliff Generate Date" timestamps on all
# three sheets, widening the "Status" column(s) to fit the new text.

$wb = $excel.ActiveWorkbook

# Excel's ColumnWidth setter snaps to whole-pixel boundaries (pixels =
# round(ColumnWidth*6+5), stored width = pixels/6), so the nearest
# reachable stored width to the target 17.2159881591797 is 17.1666...,
# produced by any ColumnWidth in [16.0834, 16.4167). Use 16.3.
$newWidth = 16.3

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-31 03:06:54"
$wsOverview.Range("E:E").ColumnWidth = $newWidth
$wsOverview.Range("F:F").ColumnWidth = $newWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-31 03:06:50"
$wsZhCn.Range("C:C").ColumnWidth = $newWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-31 03:06:54"
$wsDeDe.Range("C:C").ColumnWidth = $newWidth
